# Update gh-pages data output (generated at 456a3b4):
# refresh "want-to-go" head counts (column F) and a couple of lowest-price
# figures (column G) across the 展览 / 演出 / 本地生活 / 全部类型 sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 331
$ws.Range("F3").Value = 1148
$ws.Range("F6").Value = 3399
$ws.Range("G6").Value = 89
$ws.Range("F7").Value = 62
$ws.Range("F8").Value = 1183
$ws.Range("F9").Value = 777
$ws.Range("F10").Value = 602
$ws.Range("F12").Value = 160
$ws.Range("F13").Value = 655
$ws.Range("F14").Value = 1834
$ws.Range("F15").Value = 56
$ws.Range("F16").Value = 413
$ws.Range("F17").Value = 60
$ws.Range("F19").Value = 691
$ws.Range("F20").Value = 459
$ws.Range("F22").Value = 813
$ws.Range("F23").Value = 80269
$ws.Range("F24").Value = 80270
$ws.Range("F26").Value = 676
$ws.Range("F27").Value = 33918
$ws.Range("F28").Value = 33918
$ws.Range("F29").Value = 542
$ws.Range("F30").Value = 31
$ws.Range("F31").Value = 27
$ws.Range("F32").Value = 60
$ws.Range("F33").Value = 53
$ws.Range("F34").Value = 1005
$ws.Range("F35").Value = 316
$ws.Range("F36").Value = 163
$ws.Range("F37").Value = 639
$ws.Range("F38").Value = 2832
$ws.Range("F39").Value = 2832
$ws.Range("F40").Value = 1221
$ws.Range("F41").Value = 5511
$ws.Range("F42").Value = 800
$ws.Range("F43").Value = 460
$ws.Range("F47").Value = 434
$ws.Range("F49").Value = 7
$ws.Range("F51").Value = 58
$ws.Range("F52").Value = 8

$ws = $wb.Worksheets.Item("演出")
$ws.Range("G4").Value = "不可售"
$ws.Range("F9").Value = 1795
$ws.Range("F11").Value = 1985
$ws.Range("F12").Value = 33
$ws.Range("F14").Value = 85
$ws.Range("F15").Value = 419
$ws.Range("F18").Value = 78
$ws.Range("F20").Value = 540
$ws.Range("F32").Value = 1671
$ws.Range("F33").Value = 499
$ws.Range("F37").Value = 117
$ws.Range("F38").Value = 117
$ws.Range("F42").Value = 36
$ws.Range("F45").Value = 829
$ws.Range("F46").Value = 218
$ws.Range("F47").Value = 48

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F5").Value = 586
$ws.Range("F6").Value = 614
$ws.Range("F7").Value = 168

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 331
$ws.Range("F4").Value = 1148
$ws.Range("F6").Value = 3399
$ws.Range("G6").Value = 89
$ws.Range("F7").Value = 1183
$ws.Range("F8").Value = 777
$ws.Range("F9").Value = 614
$ws.Range("F10").Value = 614
$ws.Range("F12").Value = 1795
$ws.Range("F13").Value = 602
$ws.Range("F15").Value = 160
$ws.Range("F16").Value = 655
$ws.Range("F17").Value = 168
$ws.Range("F18").Value = 1834
$ws.Range("F19").Value = 33
$ws.Range("F20").Value = 56
$ws.Range("F21").Value = 413
$ws.Range("F22").Value = 60
$ws.Range("F24").Value = 813
$ws.Range("F27").Value = 80270
$ws.Range("F28").Value = 676
$ws.Range("F29").Value = 33918
$ws.Range("F30").Value = 542
$ws.Range("F31").Value = 31
$ws.Range("F32").Value = 27
$ws.Range("F33").Value = 540
$ws.Range("F34").Value = 540
$ws.Range("F35").Value = 53
$ws.Range("F38").Value = 316
$ws.Range("F39").Value = 163
$ws.Range("F41").Value = 2832
$ws.Range("F42").Value = 1221
$ws.Range("F43").Value = 5511
$ws.Range("F44").Value = 800
$ws.Range("F45").Value = 1671
$ws.Range("F47").Value = 117
$ws.Range("F49").Value = 434
$ws.Range("F50").Value = 36
$ws.Range("F52").Value = 7
$ws.Range("F53").Value = 218
$ws.Range("F54").Value = 58
$ws.Range("F55").Value = 8
